# Reorder slides: move the "Test Driven Development" slide (currently slide 4)
# so that it comes right before the "Structure" slide (currently slide 3).
# This swaps the positions of the two slides, matching the template-pattern
# reordering described in the commit message.

$p = $ppt.ActivePresentation

$structureIndex = 0
$tddIndex = 0

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
    if ($title -eq "Structure") {
        $structureIndex = $i
    }
    if ($title -eq "Test Driven Development") {
        $tddIndex = $i
    }
}

if ($tddIndex -gt $structureIndex) {
    $p.Slides.Item($tddIndex).MoveTo($structureIndex)
}
else {
    $p.Slides.Item($structureIndex).MoveTo($tddIndex)
}
